$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from E1 (the "c" header) into the new F1 header cell
# before renaming E1, so F1 ends up with the same style as E1.
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

# Rename existing header "c" (E1) to "cR"
$ws.Range("E1").Value = "cR"

# Set new header "cM" in F1
$ws.Range("F1").Value = "cM"

# Fill column F values (rows 2-7)
$ws.Range("F2").Value = 6
$ws.Range("F3").Value = 6
$ws.Range("F4").Value = 1
$ws.Range("F5").Value = 15
$ws.Range("F6").Value = 6
$ws.Range("F7").Value = 6

# Update the active selection to F1
$ws.Range("F1").Select()
